# Apply the edits described by the diff:
#  - Widen column A on "Human Resources" sheet
#  - Add column A width on "Marketing" sheet (previously default/no <cols>)
#  - Add column A & B widths on "Helpdesk" sheet (previously default/no <cols>)
#  - Make "Helpdesk" the active/selected sheet tab (was "Human Resources")
#  - Change the selection on "Helpdesk" sheet from D2 to C6

$wb = $excel.ActiveWorkbook

$wsHR = $wb.Worksheets.Item("Human Resources")
$wsMarketing = $wb.Worksheets.Item("Marketing")
$wsHelpdesk = $wb.Worksheets.Item("Helpdesk")

# --- Column width changes ---
# Target stored widths (OOXML <col width=".."/>): 32.42578125 on HR!A,
# 36 on Marketing!A, 35 on Helpdesk!A, 14.5703125 on Helpdesk!B.
# The ColumnWidth COM property is offset from the stored width by 5/6,
# so subtract 5/6 from each desired stored width.
$wsHR.Columns.Item(1).ColumnWidth = 32.42578125 - 0.8333333333333334
$wsMarketing.Columns.Item(1).ColumnWidth = 36 - 0.8333333333333334
$wsHelpdesk.Columns.Item(1).ColumnWidth = 35 - 0.8333333333333334
$wsHelpdesk.Columns.Item(2).ColumnWidth = 14.5703125 - 0.8333333333333334

# --- Selection on Helpdesk sheet moves from D2 to C6 ---
$wsHelpdesk.Range("C6").Select()

# --- Make Helpdesk the active tab (also clears tabSelected on Human Resources) ---
$wsHelpdesk.Activate()
